# Daily Update 키워드 10개
# The source listing re-ranked several products: the non-rank columns
# (title/link/image/price/mall/brand/manufacturer/category) of each row in
# the groups below cyclically shift to a neighboring row while column A
# (the fixed rank index) stays put.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 22, 23, 24 (source -> destination shift)
# row 22 <- former row 23 content
$ws.Cells.Item(22, 2).Value = '아이레 게이밍조립컴퓨터 배틀그라운드 오버워치 리그오브레전드 롤 본체 게임용 게이밍PC'
$ws.Cells.Item(22, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=82427199759'
$ws.Cells.Item(22, 4).Value = 'https://shopping-phinf.pstatic.net/main_8242719/82427199759.3.jpg'
$ws.Cells.Item(22, 5).Value = "'599000"
$ws.Cells.Item(22, 6).Value = ''
$ws.Cells.Item(22, 7).Value = '아이레코리아'
$ws.Cells.Item(22, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(22, 9).Value = ''
$ws.Cells.Item(22, 10).Value = ''
$ws.Cells.Item(22, 11).Value = '디지털/가전'
$ws.Cells.Item(22, 12).Value = 'PC'
$ws.Cells.Item(22, 13).Value = '조립/베어본PC'
$ws.Cells.Item(22, 14).Value = ''

# row 23 <- former row 24 content
$ws.Cells.Item(23, 2).Value = '프리미엄 게이밍 조립컴퓨터 배틀그라운드 오버워치 롤 배그 배틀필드 컴퓨터본체 견적'
$ws.Cells.Item(23, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=80277928794'
$ws.Cells.Item(23, 4).Value = 'https://shopping-phinf.pstatic.net/main_8027792/80277928794.21.jpg'
$ws.Cells.Item(23, 5).Value = "'699000"
$ws.Cells.Item(23, 6).Value = ''
$ws.Cells.Item(23, 7).Value = '쿨젠컴퓨터'
$ws.Cells.Item(23, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(23, 9).Value = '쿨젠컴퓨터'
$ws.Cells.Item(23, 10).Value = '쿨젠'
$ws.Cells.Item(23, 11).Value = '디지털/가전'
$ws.Cells.Item(23, 12).Value = 'PC'
$ws.Cells.Item(23, 13).Value = '조립/베어본PC'
$ws.Cells.Item(23, 14).Value = ''

# row 24 <- former row 22 content
$ws.Cells.Item(24, 2).Value = '라이젠 5800X GTX1650 RTX3060 Ti 3070 게이밍컴퓨터 영상편집 배그'
$ws.Cells.Item(24, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=83545548964'
$ws.Cells.Item(24, 4).Value = 'https://shopping-phinf.pstatic.net/main_8354554/83545548964.jpg'
$ws.Cells.Item(24, 5).Value = "'1375000"
$ws.Cells.Item(24, 6).Value = ''
$ws.Cells.Item(24, 7).Value = '블루컴퓨터 BLUECOM'
$ws.Cells.Item(24, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(24, 9).Value = ''
$ws.Cells.Item(24, 10).Value = ''
$ws.Cells.Item(24, 11).Value = '디지털/가전'
$ws.Cells.Item(24, 12).Value = 'PC'
$ws.Cells.Item(24, 13).Value = '조립/베어본PC'
$ws.Cells.Item(24, 14).Value = ''

# Rows 25, 26 (source -> destination shift)
# row 25 <- former row 26 content
$ws.Cells.Item(25, 2).Value = 'LG데스크탑 Z50/70 i5-3470 SSD장착 사무용컴퓨터본체'
$ws.Cells.Item(25, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=18796925108'
$ws.Cells.Item(25, 4).Value = 'https://shopping-phinf.pstatic.net/main_1879692/18796925108.2.jpg'
$ws.Cells.Item(25, 5).Value = "'235000"
$ws.Cells.Item(25, 6).Value = ''
$ws.Cells.Item(25, 7).Value = '인터파크'
$ws.Cells.Item(25, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(25, 9).Value = 'LG전자'
$ws.Cells.Item(25, 10).Value = 'LG전자'
$ws.Cells.Item(25, 11).Value = '디지털/가전'
$ws.Cells.Item(25, 12).Value = 'PC'
$ws.Cells.Item(25, 13).Value = '조립/베어본PC'
$ws.Cells.Item(25, 14).Value = ''

# row 26 <- former row 25 content
$ws.Cells.Item(26, 2).Value = '어도비 포토샵 컴퓨터 일러스트 디자인용 PC 오토캐드 본체 사양'
$ws.Cells.Item(26, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=82148907581'
$ws.Cells.Item(26, 4).Value = 'https://shopping-phinf.pstatic.net/main_8214890/82148907581.6.jpg'
$ws.Cells.Item(26, 5).Value = "'490000"
$ws.Cells.Item(26, 6).Value = ''
$ws.Cells.Item(26, 7).Value = '메가컴퓨터'
$ws.Cells.Item(26, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(26, 9).Value = ''
$ws.Cells.Item(26, 10).Value = ''
$ws.Cells.Item(26, 11).Value = '디지털/가전'
$ws.Cells.Item(26, 12).Value = 'PC'
$ws.Cells.Item(26, 13).Value = '조립/베어본PC'
$ws.Cells.Item(26, 14).Value = ''

# Rows 36, 37 (source -> destination shift)
# row 36 <- former row 37 content
$ws.Cells.Item(36, 2).Value = '게이밍 컴퓨터 조립PC RTX2060 컴퓨터본체 견적 윈도우 포함'
$ws.Cells.Item(36, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=83560670111'
$ws.Cells.Item(36, 4).Value = 'https://shopping-phinf.pstatic.net/main_8356067/83560670111.jpg'
$ws.Cells.Item(36, 5).Value = "'648640"
$ws.Cells.Item(36, 6).Value = ''
$ws.Cells.Item(36, 7).Value = '인터클러스터시스템'
$ws.Cells.Item(36, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(36, 9).Value = ''
$ws.Cells.Item(36, 10).Value = ''
$ws.Cells.Item(36, 11).Value = '디지털/가전'
$ws.Cells.Item(36, 12).Value = 'PC'
$ws.Cells.Item(36, 13).Value = '조립/베어본PC'
$ws.Cells.Item(36, 14).Value = ''

# row 37 <- former row 36 content
$ws.Cells.Item(37, 2).Value = '[하이마트]포유컴퓨터 라이젠 R5 컴퓨터본체(5600X/RX6600XT)조립PC'
$ws.Cells.Item(37, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=28544195291'
$ws.Cells.Item(37, 4).Value = 'https://shopping-phinf.pstatic.net/main_2854419/28544195291.jpg'
$ws.Cells.Item(37, 5).Value = "'1556000"
$ws.Cells.Item(37, 6).Value = ''
$ws.Cells.Item(37, 7).Value = '하이마트쇼핑몰'
$ws.Cells.Item(37, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(37, 9).Value = ''
$ws.Cells.Item(37, 10).Value = '포유컴퓨터'
$ws.Cells.Item(37, 11).Value = '디지털/가전'
$ws.Cells.Item(37, 12).Value = 'PC'
$ws.Cells.Item(37, 13).Value = '조립/베어본PC'
$ws.Cells.Item(37, 14).Value = ''

# Rows 44, 45 (source -> destination shift)
# row 44 <- former row 45 content
$ws.Cells.Item(44, 2).Value = '삼성 조립 컴퓨터 본체 27인치 모니터 세트 롤 메이플 던파'
$ws.Cells.Item(44, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=82969642005'
$ws.Cells.Item(44, 4).Value = 'https://shopping-phinf.pstatic.net/main_8296964/82969642005.jpg'
$ws.Cells.Item(44, 5).Value = "'589000"
$ws.Cells.Item(44, 6).Value = ''
$ws.Cells.Item(44, 7).Value = '윤자씨네'
$ws.Cells.Item(44, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(44, 9).Value = ''
$ws.Cells.Item(44, 10).Value = ''
$ws.Cells.Item(44, 11).Value = '디지털/가전'
$ws.Cells.Item(44, 12).Value = 'PC'
$ws.Cells.Item(44, 13).Value = '조립/베어본PC'
$ws.Cells.Item(44, 14).Value = ''

# row 45 <- former row 44 content
$ws.Cells.Item(45, 2).Value = '컴퓨터 중고 본체 조립 i5 6500 6세대 GTX1060 3G'
$ws.Cells.Item(45, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=83009330619'
$ws.Cells.Item(45, 4).Value = 'https://shopping-phinf.pstatic.net/main_8300933/83009330619.jpg'
$ws.Cells.Item(45, 5).Value = "'590000"
$ws.Cells.Item(45, 6).Value = ''
$ws.Cells.Item(45, 7).Value = '리메이드컴퓨터'
$ws.Cells.Item(45, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(45, 9).Value = ''
$ws.Cells.Item(45, 10).Value = ''
$ws.Cells.Item(45, 11).Value = '디지털/가전'
$ws.Cells.Item(45, 12).Value = 'PC'
$ws.Cells.Item(45, 13).Value = '조립/베어본PC'
$ws.Cells.Item(45, 14).Value = ''

# Rows 60, 61, 62 (source -> destination shift)
# row 60 <- former row 61 content
$ws.Cells.Item(60, 2).Value = 'i7 12700KF RTX3080 컴퓨터 본체 RB275 게이밍 PC'
$ws.Cells.Item(60, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=83630041028'
$ws.Cells.Item(60, 4).Value = 'https://shopping-phinf.pstatic.net/main_8363004/83630041028.jpg'
$ws.Cells.Item(60, 5).Value = "'4221000"
$ws.Cells.Item(60, 6).Value = ''
$ws.Cells.Item(60, 7).Value = '리메이드컴퓨터'
$ws.Cells.Item(60, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(60, 9).Value = ''
$ws.Cells.Item(60, 10).Value = ''
$ws.Cells.Item(60, 11).Value = '디지털/가전'
$ws.Cells.Item(60, 12).Value = 'PC'
$ws.Cells.Item(60, 13).Value = '조립/베어본PC'
$ws.Cells.Item(60, 14).Value = ''

# row 61 <- former row 62 content
$ws.Cells.Item(61, 2).Value = '포유컴퓨터 인텔 게이밍 i7 컴퓨터본체(10700/16G/250G)조립PC'
$ws.Cells.Item(61, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=28590326158'
$ws.Cells.Item(61, 4).Value = 'https://shopping-phinf.pstatic.net/main_2859032/28590326158.jpg'
$ws.Cells.Item(61, 5).Value = "'613000"
$ws.Cells.Item(61, 6).Value = ''
$ws.Cells.Item(61, 7).Value = '롯데ON'
$ws.Cells.Item(61, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(61, 9).Value = ''
$ws.Cells.Item(61, 10).Value = '포유컴퓨터'
$ws.Cells.Item(61, 11).Value = '디지털/가전'
$ws.Cells.Item(61, 12).Value = 'PC'
$ws.Cells.Item(61, 13).Value = '조립/베어본PC'
$ws.Cells.Item(61, 14).Value = ''

# row 62 <- former row 60 content
$ws.Cells.Item(62, 2).Value = '10세대 게이밍컴퓨터 배틀그라운드 리니지W 로스트아크 디아2레저렉션 PC 본체'
$ws.Cells.Item(62, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=82145330176'
$ws.Cells.Item(62, 4).Value = 'https://shopping-phinf.pstatic.net/main_8214533/82145330176.8.jpg'
$ws.Cells.Item(62, 5).Value = "'658000"
$ws.Cells.Item(62, 6).Value = ''
$ws.Cells.Item(62, 7).Value = '아이디컴퓨터'
$ws.Cells.Item(62, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(62, 9).Value = ''
$ws.Cells.Item(62, 10).Value = '아이디컴퓨터'
$ws.Cells.Item(62, 11).Value = '디지털/가전'
$ws.Cells.Item(62, 12).Value = 'PC'
$ws.Cells.Item(62, 13).Value = '조립/베어본PC'
$ws.Cells.Item(62, 14).Value = ''

# Rows 79, 80, 81 (source -> destination shift)
# row 79 <- former row 80 content
$ws.Cells.Item(79, 2).Value = '중고 컴퓨터 본체 i7 7700 7세대 GTX1070 8G'
$ws.Cells.Item(79, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=83009330694'
$ws.Cells.Item(79, 4).Value = 'https://shopping-phinf.pstatic.net/main_8300933/83009330694.jpg'
$ws.Cells.Item(79, 5).Value = "'1126000"
$ws.Cells.Item(79, 6).Value = ''
$ws.Cells.Item(79, 7).Value = '리메이드컴퓨터'
$ws.Cells.Item(79, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(79, 9).Value = ''
$ws.Cells.Item(79, 10).Value = ''
$ws.Cells.Item(79, 11).Value = '디지털/가전'
$ws.Cells.Item(79, 12).Value = 'PC'
$ws.Cells.Item(79, 13).Value = '조립/베어본PC'
$ws.Cells.Item(79, 14).Value = ''

# row 80 <- former row 81 content
$ws.Cells.Item(80, 2).Value = 'i7 11700F RTX3060Ti 배그전용 방송용 고사양 컴퓨터 본체PC'
$ws.Cells.Item(80, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=83600121437'
$ws.Cells.Item(80, 4).Value = 'https://shopping-phinf.pstatic.net/main_8360012/83600121437.jpg'
$ws.Cells.Item(80, 5).Value = "'2147000"
$ws.Cells.Item(80, 6).Value = ''
$ws.Cells.Item(80, 7).Value = '컴스컴바인'
$ws.Cells.Item(80, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(80, 9).Value = ''
$ws.Cells.Item(80, 10).Value = ''
$ws.Cells.Item(80, 11).Value = '디지털/가전'
$ws.Cells.Item(80, 12).Value = 'PC'
$ws.Cells.Item(80, 13).Value = '조립/베어본PC'
$ws.Cells.Item(80, 14).Value = ''

# row 81 <- former row 79 content
$ws.Cells.Item(81, 2).Value = '조립컴퓨터 AMD 풀세트 배틀그라운드 롤 RX6600 사무용 본체 견적 디아2 레저렉션'
$ws.Cells.Item(81, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=83560905478'
$ws.Cells.Item(81, 4).Value = 'https://shopping-phinf.pstatic.net/main_8356090/83560905478.2.jpg'
$ws.Cells.Item(81, 5).Value = "'599000"
$ws.Cells.Item(81, 6).Value = ''
$ws.Cells.Item(81, 7).Value = '초이스컴'
$ws.Cells.Item(81, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(81, 9).Value = '초이스컴'
$ws.Cells.Item(81, 10).Value = '초이스컴'
$ws.Cells.Item(81, 11).Value = '디지털/가전'
$ws.Cells.Item(81, 12).Value = 'PC'
$ws.Cells.Item(81, 13).Value = '조립/베어본PC'
$ws.Cells.Item(81, 14).Value = ''

# Rows 82, 83 (source -> destination shift)
# row 82 <- former row 83 content
$ws.Cells.Item(82, 2).Value = '피씨홀릭 배그용컴퓨터 조립PC 컴퓨터 본체 라이젠 3500x 2060 super'
$ws.Cells.Item(82, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=82686440531'
$ws.Cells.Item(82, 4).Value = 'https://shopping-phinf.pstatic.net/main_8268644/82686440531.jpg'
$ws.Cells.Item(82, 5).Value = "'1599000"
$ws.Cells.Item(82, 6).Value = ''
$ws.Cells.Item(82, 7).Value = 'PC Holic'
$ws.Cells.Item(82, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(82, 9).Value = ''
$ws.Cells.Item(82, 10).Value = ''
$ws.Cells.Item(82, 11).Value = '디지털/가전'
$ws.Cells.Item(82, 12).Value = 'PC'
$ws.Cells.Item(82, 13).Value = '조립/베어본PC'
$ws.Cells.Item(82, 14).Value = ''

# row 83 <- former row 82 content
$ws.Cells.Item(83, 2).Value = '인텔12세대 i5 컴퓨터본체(12400F/RX6600)조립PC'
$ws.Cells.Item(83, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=29621216094'
$ws.Cells.Item(83, 4).Value = 'https://shopping-phinf.pstatic.net/main_2962121/29621216094.jpg'
$ws.Cells.Item(83, 5).Value = "'1369000"
$ws.Cells.Item(83, 6).Value = ''
$ws.Cells.Item(83, 7).Value = '인터파크'
$ws.Cells.Item(83, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(83, 9).Value = ''
$ws.Cells.Item(83, 10).Value = '포유컴퓨터'
$ws.Cells.Item(83, 11).Value = '디지털/가전'
$ws.Cells.Item(83, 12).Value = 'PC'
$ws.Cells.Item(83, 13).Value = '조립/베어본PC'
$ws.Cells.Item(83, 14).Value = ''

# Rows 90, 91 (source -> destination shift)
# row 90 <- former row 91 content
$ws.Cells.Item(90, 2).Value = 'LG S급 외관 슬림형 초고속 i3 i5 사무용 컴퓨터 본체'
$ws.Cells.Item(90, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=18859911510'
$ws.Cells.Item(90, 4).Value = 'https://shopping-phinf.pstatic.net/main_1885991/18859911510.1.jpg'
$ws.Cells.Item(90, 5).Value = "'198550"
$ws.Cells.Item(90, 6).Value = ''
$ws.Cells.Item(90, 7).Value = '인터파크'
$ws.Cells.Item(90, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(90, 9).Value = 'LG전자'
$ws.Cells.Item(90, 10).Value = 'LG전자'
$ws.Cells.Item(90, 11).Value = '디지털/가전'
$ws.Cells.Item(90, 12).Value = 'PC'
$ws.Cells.Item(90, 13).Value = '조립/베어본PC'
$ws.Cells.Item(90, 14).Value = ''

# row 91 <- former row 90 content
$ws.Cells.Item(91, 2).Value = '완본체 영상편집용 배그용 중고컴퓨터본체 서든컴퓨터'
$ws.Cells.Item(91, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=83705242553'
$ws.Cells.Item(91, 4).Value = 'https://shopping-phinf.pstatic.net/main_8370524/83705242553.jpg'
$ws.Cells.Item(91, 5).Value = "'1199000"
$ws.Cells.Item(91, 6).Value = ''
$ws.Cells.Item(91, 7).Value = '아이비컴즈'
$ws.Cells.Item(91, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(91, 9).Value = ''
$ws.Cells.Item(91, 10).Value = ''
$ws.Cells.Item(91, 11).Value = '디지털/가전'
$ws.Cells.Item(91, 12).Value = 'PC'
$ws.Cells.Item(91, 13).Value = '조립/베어본PC'
$ws.Cells.Item(91, 14).Value = ''

# Rows 92, 93, 95, 94 (source -> destination shift)
# row 92 <- former row 93 content
$ws.Cells.Item(92, 2).Value = '라이젠 게이밍 R5 컴퓨터본체(5600X/GTX1650)조립PC'
$ws.Cells.Item(92, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=82254931550'
$ws.Cells.Item(92, 4).Value = 'https://shopping-phinf.pstatic.net/main_8225493/82254931550.jpg'
$ws.Cells.Item(92, 5).Value = "'1036000"
$ws.Cells.Item(92, 6).Value = ''
$ws.Cells.Item(92, 7).Value = '주식회사 포유컴퓨터'
$ws.Cells.Item(92, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(92, 9).Value = 'FORYOUCOM'
$ws.Cells.Item(92, 10).Value = ''
$ws.Cells.Item(92, 11).Value = '디지털/가전'
$ws.Cells.Item(92, 12).Value = 'PC'
$ws.Cells.Item(92, 13).Value = '조립/베어본PC'
$ws.Cells.Item(92, 14).Value = ''

# row 93 <- former row 95 content
$ws.Cells.Item(93, 2).Value = 'i7 12700KF RTX3070 8GB 컴퓨터 본체 RB273 게이밍 PC'
$ws.Cells.Item(93, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=83630095060'
$ws.Cells.Item(93, 4).Value = 'https://shopping-phinf.pstatic.net/main_8363009/83630095060.jpg'
$ws.Cells.Item(93, 5).Value = "'3247000"
$ws.Cells.Item(93, 6).Value = ''
$ws.Cells.Item(93, 7).Value = '리메이드컴퓨터'
$ws.Cells.Item(93, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(93, 9).Value = ''
$ws.Cells.Item(93, 10).Value = ''
$ws.Cells.Item(93, 11).Value = '디지털/가전'
$ws.Cells.Item(93, 12).Value = 'PC'
$ws.Cells.Item(93, 13).Value = '조립/베어본PC'
$ws.Cells.Item(93, 14).Value = ''

# row 95 <- former row 94 content
$ws.Cells.Item(95, 2).Value = '프리플로우 GAMING R5 R6 컴퓨터본체,AMD 라이젠5 5600X,라데온 RX6600 조립PC'
$ws.Cells.Item(95, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=29224410733'
$ws.Cells.Item(95, 4).Value = 'https://shopping-phinf.pstatic.net/main_2922441/29224410733.1.jpg'
$ws.Cells.Item(95, 5).Value = "'1479000"
$ws.Cells.Item(95, 6).Value = ''
$ws.Cells.Item(95, 7).Value = '롯데ON'
$ws.Cells.Item(95, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(95, 9).Value = '프리플로우'
$ws.Cells.Item(95, 10).Value = '프리플로우'
$ws.Cells.Item(95, 11).Value = '디지털/가전'
$ws.Cells.Item(95, 12).Value = 'PC'
$ws.Cells.Item(95, 13).Value = '조립/베어본PC'
$ws.Cells.Item(95, 14).Value = ''

# row 94 <- former row 92 content
$ws.Cells.Item(94, 2).Value = '삼성전자 화이트 게이밍 조립식 i5-7500 GTX1060 컴퓨터 본체'
$ws.Cells.Item(94, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=26580822837'
$ws.Cells.Item(94, 4).Value = 'https://shopping-phinf.pstatic.net/main_2658082/26580822837.jpg'
$ws.Cells.Item(94, 5).Value = "'688500"
$ws.Cells.Item(94, 6).Value = ''
$ws.Cells.Item(94, 7).Value = 'G마켓'
$ws.Cells.Item(94, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(94, 9).Value = '삼성'
$ws.Cells.Item(94, 10).Value = '삼성전자'
$ws.Cells.Item(94, 11).Value = '디지털/가전'
$ws.Cells.Item(94, 12).Value = 'PC'
$ws.Cells.Item(94, 13).Value = '조립/베어본PC'
$ws.Cells.Item(94, 14).Value = ''

# Rows 99, 100 (source -> destination shift)
# row 99 <- former row 100 content
$ws.Cells.Item(99, 2).Value = '인텔 게이밍 i7 컴퓨터본체(10700/16G/250G)조립PC'
$ws.Cells.Item(99, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=7958233834'
$ws.Cells.Item(99, 4).Value = 'https://shopping-phinf.pstatic.net/main_7958233/7958233834.jpg'
$ws.Cells.Item(99, 5).Value = "'653000"
$ws.Cells.Item(99, 6).Value = ''
$ws.Cells.Item(99, 7).Value = '주식회사 포유컴퓨터'
$ws.Cells.Item(99, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(99, 9).Value = 'FORYOUCOM'
$ws.Cells.Item(99, 10).Value = ''
$ws.Cells.Item(99, 11).Value = '디지털/가전'
$ws.Cells.Item(99, 12).Value = 'PC'
$ws.Cells.Item(99, 13).Value = '조립/베어본PC'
$ws.Cells.Item(99, 14).Value = ''

# row 100 <- former row 99 content
$ws.Cells.Item(100, 2).Value = '12400F 12세대 RTX3060 PC 컴퓨터 게이밍 본체 RB242'
$ws.Cells.Item(100, 3).Value = 'https://search.shopping.naver.com/gate.nhn?id=83706280315'
$ws.Cells.Item(100, 4).Value = 'https://shopping-phinf.pstatic.net/main_8370628/83706280315.jpg'
$ws.Cells.Item(100, 5).Value = "'1662000"
$ws.Cells.Item(100, 6).Value = ''
$ws.Cells.Item(100, 7).Value = '리메이드컴퓨터'
$ws.Cells.Item(100, 8).Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Cells.Item(100, 9).Value = ''
$ws.Cells.Item(100, 10).Value = ''
$ws.Cells.Item(100, 11).Value = '디지털/가전'
$ws.Cells.Item(100, 12).Value = 'PC'
$ws.Cells.Item(100, 13).Value = '조립/베어본PC'
$ws.Cells.Item(100, 14).Value = ''

